$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "60.118.66"
$ws.Range("E2").Value = "  +0.23%  "

$ws.Range("D3").Value = "2.418.76"
$ws.Range("E3").Value = "  +0.11%  "

$ws.Range("E4").Value = "  -0.06%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "553.65"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +0.30%  "

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "136.85"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  -0.01%  "

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.587"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  +0.92%  "

$ws.Range("E9").Value = "  -0.70%  "

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "5.69"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  -1.84%  "

$ws.Range("E11").Value = "  -0.48%  "

$ws.Range("E12").Value = "  -1.03%  "

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "24.81"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  -0.12%  "

$ws.Range("D14").Value = "2.849.23"
$ws.Range("E14").Value = "  -0.02%  "

$ws.Range("D15").Value = "60.010.84"
$ws.Range("E15").Value = "  +0.12%  "

$ws.Range("E16").Value = "  -0.22%  "

$ws.Range("D17").Value = "2.432.69"
$ws.Range("E17").Value = "  +1.11%  "

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "11.25"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  -1.23%  "

$ws.Range("E19").Value = "  +2.66%  "

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "327.76"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  -1.06%  "

$ws.Range("E21").Value = "  +1.16%  "

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "1.00"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  +0.19%  "

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "64.82"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  -1.14%  "

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "0.179"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  +4.91%  "

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "8.62"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  +0.14%  "

$ws.Range("E26").Value = "  +0.10%  "

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "1.42"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  +6.17%  "

$ws.Range("D28").Value = "0.0₃0773"
$ws.Range("E28").Value = "  -0.95%  "

$ws.Range("E29").Value = "  -0.08%  "

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "170.54"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  +0.02%  "

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "6.11"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  -1.74%  "

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "1.07"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  +4.13%  "

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "0.402"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  -3.13%  "

$ws.Range("E34").Value = "  -0.40%  "

$ws.Range("E35").Value = "  +0.04%  "

$ws.Range("E36").Value = "  +2.77%  "

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "4.23"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  +1.08%  "

$ws.Range("E38").Value = "  +0.02%  "

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "322.47"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  +3.22%  "

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "1.60"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  -1.01%  "

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "145.88"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  +4.56%  "

$ws.Range("E42").Value = "  -0.80%  "

$ws.Range("E43").Value = "  +0.17%  "

$ws.Range("E44").Value = "  +2.83%  "

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.0516"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  -0.79%  "

$ws.Range("E46").Value = "  +0.60%  "

$ws.Range("E47").Value = "  -1.15%  "

$ws.Range("E48").Value = "  -0.15%  "

$ws.Range("E49").Value = "  -0.73%  "

$ws.Range("E50").Value = "  -0.74%  "

$ws.Range("E51").Value = "  -1.53%  "

